$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Judge1 (G) / Judge2 (H) for existing rows 2-49 ---
$judges = @(
    @('Yaoying Wu', 'Era Jain'),
    @('Anupam Pandey', 'Sucheta Soundarajan'),
    @('Mary Beth Monroe', 'Sucheta Soundarajan'),
    @('Yaoying Wu', 'Baris Salman'),
    @('Bing Dong', 'Farzana Rahman'),
    @('M. Cenk Gursoy', 'Yiyang Sun'),
    @('Shikha Nangia', 'Jay Henderson'),
    @('Era Jain', 'Farzana Rahman'),
    @('Wanliang Shan', 'Qinru Qiu'),
    @('Ben Akih-Kumgeh', 'Anupam Pandey'),
    @('Svetoslava Todorova', 'Baris Salman'),
    @('Shikha Nangia', 'Yaoying Wu'),
    @('Yi Zheng', 'Mary Beth Monroe'),
    @('Nadeem Ghani', 'Jesse Q. Bond'),
    @('Ashok Sangani', 'Senem Velipasalar'),
    @('Ruth Chen', 'Bing Dong'),
    @('Ruth Chen', 'Jeongmin Ahn'),
    @('Chikukuri Monhan', 'Ashok Sangani'),
    @('Yaoying Wu', 'Jay Henderson'),
    @('Jason Pollack', 'Zhenyu Gan'),
    @('Anupam Pandey', 'Yuzhe Tang'),
    @('Svetoslava Todorova', 'Elizabeth Carter'),
    @('Gabriel Silva De Oliveira', 'Mary Beth Monroe'),
    @('Jesse Q. Bond', 'Sucheta Soundarajan'),
    @('Ian Hosein', 'Endadul Hoque'),
    @('Amit Sanyal', 'Yuzhe Tang'),
    @('Ben Akih-Kumgeh', 'Joao Paulo Marum'),
    @('Anupam Pandey', 'Gabriel Silva De Oliveira'),
    @('Jesse Q. Bond', 'Shikha Nangia'),
    @('Joao Paulo Marum', 'Elizabeth Carter'),
    @('Nadeem Ghani', 'Shalabh Maroo'),
    @('Amit Sanyal', 'Nadeem Ghani'),
    @('Chikukuri Monhan', 'Senem Velipasalar'),
    @('Yi Zheng', 'Elizabeth Carter'),
    @('Ian Hosein', 'Min Liu'),
    @('John F. Dannenhoffer', 'M. Cenk Gursoy'),
    @('Baris Salman', 'Pankaj Jha'),
    @('Yiyang Sun', 'Yi Zheng'),
    @('Ian Hosein', 'Zhenyu Gan'),
    @('Younes Radi', 'Mary Beth Monroe'),
    @('Endadul Hoque', 'Wanliang Shan'),
    @('Jeongmin Ahn', 'Ruth Chen'),
    @('Senem Velipasalar', 'Ruth Chen'),
    @('Shikha Nangia', 'Wanliang Shan'),
    @('Jason Pollack', 'Pankaj Jha'),
    @('Zhenyu Gan', 'Yuzhe Tang'),
    @('Baris Salman', 'Shalabh Maroo'),
    @('Svetoslava Todorova', 'Svetoslava Todorova')
)
for ($i = 0; $i -lt $judges.Length; $i++) {
    $pair = $judges[$i]
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $pair[0]
    $ws.Cells.Item($row, 8).Value = $pair[1]
}

# --- Append new rows 50-69 (poster entries 49-68) ---
$newRows = @(
    @(49, 'AI in Mental Health Diagnosis', 'Using machine learning to detect mental health conditions.', 'Benjamin', 'Foster', 'Psychology & AI', 'Yaoying Wu', 'Jeongmin Ahn'),
    @(50, 'Smart Prosthetics and Bionics', 'AI-powered prosthetics improving mobility.', 'Madison', 'Barnes', 'Biomedical Engineering', 'C.Y. Roger Chen', 'Ruth Chen'),
    @(51, 'The Future of Electric Aircraft', 'Advancements in battery tech for aviation.', 'Ryan', 'Henderson', 'Aerospace Engineering', 'John F. Dannenhoffer', 'Ashok Sangani'),
    @(52, 'AI-Powered Smart Homes', 'Enhancing home automation with artificial intelligence.', 'Sofia', 'Bailey', 'Electrical Engineering', 'Zhenyu Gan', 'M. Cenk Gursoy'),
    @(53, 'The Role of AI in Criminal Justice', 'Examining AI’s impact on legal decision-making.', 'Elijah', 'Jenkins', 'Law & Technology', 'Ashok Sangani', 'C.Y. Roger Chen'),
    @(54, 'Sustainable Space Exploration', 'Using renewable tech for deep space missions.', 'Abigail', 'Rivera', 'Aerospace Engineering', 'Ashok Sangani', 'Bing Dong'),
    @(55, 'AI in Sports Performance Analysis', 'How machine learning enhances athlete training.', 'Lucas', 'Brooks', 'Sports Science', 'Elizabeth Carter', 'Shikha Nangia'),
    @(56, 'Smart Traffic Management with AI', 'Reducing congestion with predictive AI systems.', 'Chloe', 'Kelly', 'Civil Engineering', 'Anupam Pandey', 'Mary Beth Monroe'),
    @(57, 'Next-Gen Biometric Authentication', 'Security advancements in facial and fingerprint recognition.', 'Daniel', 'Cooper', 'Cybersecurity', 'Shalabh Maroo', 'Zhen Ma'),
    @(58, 'AI in Archaeology', 'Using deep learning to reconstruct ancient civilizations.', 'Grace', 'Howard', 'Anthropology & AI', 'Yi Zheng', 'Yi Zheng'),
    @(59, 'Deepfake Detection with AI', 'Preventing AI-generated misinformation.', 'Oliver', 'Ward', 'Media & Communication', 'Zhenyu Gan', 'Elizabeth Carter'),
    @(60, 'Wireless Charging for Electric Vehicles', 'Improving efficiency in wireless power transfer.', 'Aiden', 'Cox', 'Electrical Engineering', 'Chikukuri Monhan', 'Joao Paulo Marum'),
    @(61, 'AI-Powered Chatbots for Healthcare', 'Virtual assistants for medical diagnosis.', 'Lily', 'Morgan', 'Biomedical Engineering', 'Senem Velipasalar', 'Zhen Ma'),
    @(62, 'Ocean Cleanup with Robotics', 'Using AI-driven robots to remove ocean waste.', 'Emma', 'Flores', 'Environmental Science', 'Senem Velipasalar', 'Yi Zheng'),
    @(63, 'AI and Personalized Nutrition', 'Machine learning for custom diet plans.', 'Carter', 'Richardson', 'Health Informatics', 'Svetoslava Todorova', 'Qinru Qiu'),
    @(64, 'Space Tourism: Challenges and Innovations', 'The feasibility of commercial space travel.', 'Mila', 'Wood', 'Aerospace Engineering', 'Min Liu', 'Younes Radi'),
    @(65, 'AI in Smart Farming', 'AI-powered tools for crop monitoring.', 'Jack', 'Bennett', 'Agricultural Engineering', 'Yaoying Wu', 'Yiyang Sun'),
    @(66, 'AI for Earthquake Prediction', 'Machine learning models for seismic activity analysis.', 'Scarlett', 'Gray', 'Geophysics', 'Chikukuri Monhan', 'Zhen Ma'),
    @(67, 'The Future of Quantum Sensors', 'Ultra-precise measurements for navigation and science.', 'Noah', 'James', 'Quantum Engineering', 'Bing Dong', 'Ruth Chen'),
    @(68, 'AI and Emotional Recognition', 'Exploring AI’s ability to detect human emotions.', 'Avery', 'Russell', 'Artificial Intelligence', 'Baris Salman', 'Ben Akih-Kumgeh')
)
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $newRows[$i]
    $row = 50 + $i
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
}
